$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 263.26315
$ws.Range("I12").Value = 238.6875
$ws.Range("J12").Value = 394.33334
$ws.Range("K12").Value = 238.6875
$ws.Range("L12").Value = 394.33334
$ws.Range("M12").Value = -68.6875
$ws.Range("N12").Value = -734.33334
$ws.Range("H32").Value = 3468
$ws.Range("I32").Value = 2827
$ws.Range("K32").Value = 2827
$ws.Range("M32").Value = -2501
$ws.Range("H64").Value = 4521.8887
$ws.Range("I64").Value = 4283
$ws.Range("J64").Value = 4999.6665
$ws.Range("K64").Value = 4283
$ws.Range("L64").Value = 4999.6665
$ws.Range("M64").Value = -4035
$ws.Range("N64").Value = -5495.6665
$ws.Range("H67").Value = 4521.8887
$ws.Range("I67").Value = 4283
$ws.Range("J67").Value = 4999.6665
$ws.Range("K67").Value = 4283
$ws.Range("L67").Value = 4999.6665
$ws.Range("M67").Value = -3425
$ws.Range("N67").Value = -6715.6665
$ws.Range("H76").Value = 7142.5713
$ws.Range("I76").Value = 5499.5
$ws.Range("K76").Value = 5499.5
$ws.Range("M76").Value = -5184.5
$ws.Range("H79").Value = 7142.5713
$ws.Range("I79").Value = 5499.5
$ws.Range("K79").Value = 5499.5
$ws.Range("M79").Value = -4407.5
$ws.Range("H86").Value = 1055496
$ws.Range("I86").Value = 1430744.6
$ws.Range("K86").Value = 1430744.6
$ws.Range("M86").Value = -1429621.6
$ws.Range("H89").Value = 1055496
$ws.Range("I89").Value = 1430744.6
$ws.Range("K89").Value = 7153723
$ws.Range("M89").Value = -7148107
$ws.Range("H92").Value = 207.78947
$ws.Range("I92").Value = 222.625
$ws.Range("K92").Value = 222.625
$ws.Range("M92").Value = 1025.375
$ws.Range("H125").Value = 2061.625
$ws.Range("I125").Value = 1299.6
$ws.Range("K125").Value = 11696.4
$ws.Range("M125").Value = -9236.4
$ws.Range("H132").Value = 2217
$ws.Range("I132").Value = 1937.4706
$ws.Range("K132").Value = 5812.4118
$ws.Range("M132").Value = -3282.4118

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3520.8628
$ws.Range("I132").Value = 2212.8125
$ws.Range("K132").Value = 6638.4375
$ws.Range("M132").Value = -4108.4375

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H94").Value = 1190.1333
$ws.Range("I94").Value = 1099.8572
$ws.Range("J94").Value = 1400.7778
$ws.Range("K94").Value = 1099.8572
$ws.Range("L94").Value = 1400.7778
$ws.Range("M94").Value = -648.8571999999999
$ws.Range("N94").Value = -2302.7778
$ws.Range("H99").Value = 3445.4092
$ws.Range("I99").Value = 3333.1667
$ws.Range("J99").Value = 3487.5
$ws.Range("K99").Value = 3333.1667
$ws.Range("L99").Value = 3487.5
$ws.Range("M99").Value = -1835.1667
$ws.Range("N99").Value = -6483.5
$ws.Range("H105").Value = 2172.0344
$ws.Range("I105").Value = 1619.2307
$ws.Range("J105").Value = 2621.1875
$ws.Range("K105").Value = 1619.2307
$ws.Range("L105").Value = 2621.1875
$ws.Range("M105").Value = 127.7692999999999
$ws.Range("N105").Value = -6115.1875
$ws.Range("H108").Value = 104954.5
$ws.Range("J108").Value = 104954.5
$ws.Range("L108").Value = 104954.5
$ws.Range("N108").Value = -112634.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 1578.6428
$ws.Range("I132").Value = 1425.1666
$ws.Range("K132").Value = 4275.4998
$ws.Range("M132").Value = -1745.4998
$ws.Range("H134").Value = 3219
$ws.Range("I134").Value = 1857.826
$ws.Range("K134").Value = 5573.478
$ws.Range("M134").Value = -3038.478

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 84124
$ws.Range("J37").Value = 84124
$ws.Range("L37").Value = 252372
$ws.Range("N37").Value = -252596
$ws.Range("H38").Value = 308.75
$ws.Range("J38").Value = 76.2
$ws.Range("L38").Value = 228.6
$ws.Range("N38").Value = -922.6
$ws.Range("H94").Value = 4062.5
$ws.Range("J94").Value = 4062.5
$ws.Range("L94").Value = 12187.5
$ws.Range("N94").Value = -13539.5
$ws.Range("H107").Value = 596.44446
$ws.Range("I107").Value = 439.65
$ws.Range("K107").Value = 1318.95
$ws.Range("M107").Value = 601.0500000000002
$ws.Range("I114").Value = 1250
$ws.Range("K114").Value = 3750
$ws.Range("M114").Value = -496
$ws.Range("H117").Value = 182
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 182
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 546
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -7430
$ws.Range("H131").Value = 5022.491
$ws.Range("I131").Value = 3735.625
$ws.Range("K131").Value = 11206.875
$ws.Range("M131").Value = -6166.875
$ws.Range("H138").Value = 3681.6365
$ws.Range("I138").Value = 1899.6
$ws.Range("J138").Value = 5166.6665
$ws.Range("K138").Value = 5698.799999999999
$ws.Range("L138").Value = 15499.9995
$ws.Range("M138").Value = -558.7999999999993
$ws.Range("N138").Value = -25779.9995
$ws.Range("H139").Value = 2799.2
$ws.Range("I139").Value = 1614.1538
$ws.Range("K139").Value = 4842.4614
$ws.Range("M139").Value = 297.5385999999999
$ws.Range("H140").Value = 179523.94
$ws.Range("I140").Value = 203060.47
$ws.Range("K140").Value = 609181.41
$ws.Range("M140").Value = -604001.41
$ws.Range("H141").Value = 346600
$ws.Range("I141").Value = 1507450
$ws.Range("J141").Value = 14928.571
$ws.Range("K141").Value = 4522350
$ws.Range("L141").Value = 44785.713
$ws.Range("M141").Value = -4517170
$ws.Range("N141").Value = -55145.713

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7412
$ws.Range("I70").Value = 6797
$ws.Range("K70").Value = 6797
$ws.Range("M70").Value = -6527
$ws.Range("H73").Value = 7412
$ws.Range("I73").Value = 6797
$ws.Range("K73").Value = 6797
$ws.Range("M73").Value = -5861
$ws.Range("H113").Value = 3357.1428
$ws.Range("I113").Value = 2125.125
$ws.Range("K113").Value = 2125.125
$ws.Range("M113").Value = 44.875
$ws.Range("H126").Value = 5002.75
$ws.Range("I126").Value = 5002.75
$ws.Range("K126").Value = 15008.25
$ws.Range("M126").Value = -12538.25
$ws.Range("H132").Value = 21281826
$ws.Range("I132").Value = 27030152
$ws.Range("J132").Value = 13024.5
$ws.Range("K132").Value = 81090456
$ws.Range("L132").Value = 39073.5
$ws.Range("M132").Value = -81087926
$ws.Range("N132").Value = -44133.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 107502.9
$ws.Range("I7").Value = 6002.5
$ws.Range("K7").Value = 6002.5
$ws.Range("M7").Value = -5890.5
$ws.Range("H74").Value = 129000
$ws.Range("I74").Value = 125000
$ws.Range("K74").Value = 125000
$ws.Range("M74").Value = -124002
$ws.Range("H77").Value = 129000
$ws.Range("I77").Value = 125000
$ws.Range("K77").Value = 375000
$ws.Range("M77").Value = -370008
$ws.Range("H126").Value = 107502.9
$ws.Range("I126").Value = 6002.5
$ws.Range("K126").Value = 18007.5
$ws.Range("M126").Value = -15537.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 52458.5
$ws.Range("J123").Value = 52458.5
$ws.Range("L123").Value = 52458.5
$ws.Range("N123").Value = -62258.5
$ws.Range("H132").Value = 2034.68
$ws.Range("I132").Value = 1353.55
$ws.Range("K132").Value = 4060.65
$ws.Range("M132").Value = -1530.65
$ws.Range("H140").Value = 53447
$ws.Range("J140").Value = 53447
$ws.Range("L140").Value = 53447
$ws.Range("N140").Value = -63807
